$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1716.8125
$ws.Range("I98").Value = 1716.8125
$ws.Range("K98").Value = 1716.8125
$ws.Range("M98").Value = -218.8125
$ws.Range("H122").Value = 1716.8125
$ws.Range("I122").Value = 1716.8125
$ws.Range("K122").Value = 5150.4375
$ws.Range("M122").Value = -2700.4375
$ws.Range("H129").Value = 1971.125
$ws.Range("I129").Value = 154.2
$ws.Range("K129").Value = 462.6
$ws.Range("M129").Value = 4537.4
$ws.Range("H138").Value = 9502.192999999999
$ws.Range("J138").Value = 9572.356
$ws.Range("L138").Value = 28717.068
$ws.Range("N138").Value = -38997.068

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 369
$ws.Range("I4").Value = 368.5
$ws.Range("K4").Value = 368.5
$ws.Range("M4").Value = -252.5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H37").Value = 32433.572
$ws.Range("J37").Value = 33172.5
$ws.Range("L37").Value = 33172.5
$ws.Range("N37").Value = -33718.5
$ws.Range("H132").Value = 3780.8286
$ws.Range("I132").Value = 3782.8
$ws.Range("K132").Value = 11348.4
$ws.Range("M132").Value = -8818.400000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 308.77777
$ws.Range("I22").Value = 308.77777
$ws.Range("K22").Value = 308.77777
$ws.Range("M22").Value = -135.77777
$ws.Range("H74").Value = 23619
$ws.Range("I74").Value = 26492.334
$ws.Range("K74").Value = 26492.334
$ws.Range("M74").Value = -25556.334
$ws.Range("H77").Value = 23619
$ws.Range("I77").Value = 26492.334
$ws.Range("K77").Value = 79477.00199999999
$ws.Range("M77").Value = -74797.00199999999
$ws.Range("H86").Value = 4439.625
$ws.Range("I86").Value = 4253
$ws.Range("K86").Value = 4253
$ws.Range("M86").Value = -3130
$ws.Range("H89").Value = 4439.625
$ws.Range("I89").Value = 4253
$ws.Range("K89").Value = 21265
$ws.Range("M89").Value = -15649
$ws.Range("H94").Value = 2963.8
$ws.Range("I94").Value = 709
$ws.Range("J94").Value = 3527.5
$ws.Range("K94").Value = 709
$ws.Range("L94").Value = 3527.5
$ws.Range("M94").Value = -258
$ws.Range("N94").Value = -4429.5
$ws.Range("H99").Value = 1778.2273
$ws.Range("I99").Value = 1778.2273
$ws.Range("K99").Value = 1778.2273
$ws.Range("M99").Value = -280.2273
$ws.Range("H105").Value = 6659.7827
$ws.Range("I105").Value = 5997.0713
$ws.Range("K105").Value = 5997.0713
$ws.Range("M105").Value = -4250.0713
$ws.Range("H107").Value = 1843.1818
$ws.Range("I107").Value = 1843.1818
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1843.1818
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 76.81819999999993
$ws.Range("N107").ClearContents()
$ws.Range("H141").Value = 80000
$ws.Range("I141").Value = 80000
$ws.Range("K141").Value = 80000
$ws.Range("M141").Value = -74820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2845.4443
$ws.Range("J31").Value = 3463.3333
$ws.Range("L31").Value = 3463.3333
$ws.Range("N31").Value = -4053.3333
$ws.Range("H34").Value = 2845.4443
$ws.Range("J34").Value = 3463.3333
$ws.Range("L34").Value = 3463.3333
$ws.Range("N34").Value = -3867.3333
$ws.Range("H58").Value = 4774.923
$ws.Range("I58").Value = 4336.913
$ws.Range("K58").Value = 4336.913
$ws.Range("M58").Value = -4133.913
$ws.Range("H122").Value = 3076.9092
$ws.Range("I122").Value = 3059.353
$ws.Range("J122").Value = 3095.5625
$ws.Range("K122").Value = 9178.059000000001
$ws.Range("L122").Value = 9286.6875
$ws.Range("M122").Value = -6728.059000000001
$ws.Range("N122").Value = -14186.6875
$ws.Range("H132").Value = 3921.45
$ws.Range("I132").Value = 3584.647
$ws.Range("K132").Value = 10753.941
$ws.Range("M132").Value = -8223.940999999999
$ws.Range("H136").Value = 4774.923
$ws.Range("I136").Value = 4336.913
$ws.Range("K136").Value = 13010.739
$ws.Range("M136").Value = -10460.739

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5668954.5
$ws.Range("I4").Value = 6540809
$ws.Range("K4").Value = 19622427
$ws.Range("M4").Value = -19622315
$ws.Range("H113").Value = 2197.6
$ws.Range("J113").Value = 2197.6
$ws.Range("L113").Value = 6592.799999999999
$ws.Range("N113").Value = -10932.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12522.308
$ws.Range("I43").Value = 5700
$ws.Range("J43").Value = 15554.444
$ws.Range("K43").Value = 5700
$ws.Range("L43").Value = 15554.444
$ws.Range("M43").Value = -5549
$ws.Range("N43").Value = -15856.444
$ws.Range("H46").Value = 18183.334
$ws.Range("J46").Value = 21000
$ws.Range("L46").Value = 21000
$ws.Range("N46").Value = -21312
$ws.Range("H57").Value = 14705.5
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26640
$ws.Range("H80").Value = 10537.6
$ws.Range("I80").Value = 5196.5713
$ws.Range("J80").Value = 23000
$ws.Range("K80").Value = 5196.5713
$ws.Range("L80").Value = 23000
$ws.Range("M80").Value = -4198.5713
$ws.Range("N80").Value = -24996
$ws.Range("H83").Value = 10537.6
$ws.Range("I83").Value = 5196.5713
$ws.Range("J83").Value = 23000
$ws.Range("K83").Value = 25982.8565
$ws.Range("L83").Value = 115000
$ws.Range("M83").Value = -20990.8565
$ws.Range("N83").Value = -124984
$ws.Range("H102").Value = 39271.96
$ws.Range("I102").Value = 45586.477
$ws.Range("K102").Value = 45586.477
$ws.Range("M102").Value = -43964.477
$ws.Range("H122").Value = 1986.5
$ws.Range("I122").Value = 1412.9
$ws.Range("K122").Value = 4238.700000000001
$ws.Range("M122").Value = -1788.700000000001
$ws.Range("H132").Value = 5366.7646
$ws.Range("I132").Value = 5568.7393
$ws.Range("J132").Value = 4944.4546
$ws.Range("K132").Value = 16706.2179
$ws.Range("L132").Value = 14833.3638
$ws.Range("M132").Value = -14176.2179
$ws.Range("N132").Value = -19893.3638

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2624.8
$ws.Range("J46").Value = 3249.8333
$ws.Range("L46").Value = 3249.8333
$ws.Range("N46").Value = -3625.8333
$ws.Range("H61").Value = 1267.7646
$ws.Range("I61").Value = 1183.2858
$ws.Range("K61").Value = 1183.2858
$ws.Range("M61").Value = -981.2858000000001
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H113").Value = 1267.7646
$ws.Range("I113").Value = 1183.2858
$ws.Range("K113").Value = 1183.2858
$ws.Range("M113").Value = 986.7141999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 9939.615
$ws.Range("I14").Value = 26512.25
$ws.Range("J14").Value = 2574
$ws.Range("K14").Value = 26512.25
$ws.Range("L14").Value = 2574
$ws.Range("M14").Value = -26344.25
$ws.Range("N14").Value = -2910
$ws.Range("H62").Value = 28890.908
$ws.Range("I62").Value = 12966.667
$ws.Range("J62").Value = 34862.5
$ws.Range("K62").Value = 12966.667
$ws.Range("L62").Value = 34862.5
$ws.Range("M62").Value = -12342.667
$ws.Range("N62").Value = -36110.5
$ws.Range("H65").Value = 28890.908
$ws.Range("I65").Value = 12966.667
$ws.Range("J65").Value = 34862.5
$ws.Range("K65").Value = 64833.335
$ws.Range("L65").Value = 174312.5
$ws.Range("M65").Value = -61713.335
$ws.Range("N65").Value = -180552.5
$ws.Range("H81").Value = 2542.4285
$ws.Range("I81").Value = 2416.1667
$ws.Range("K81").Value = 4832.3334
$ws.Range("M81").Value = -3771.3334
$ws.Range("H84").Value = 2542.4285
$ws.Range("I84").Value = 2416.1667
$ws.Range("K84").Value = 24161.667
$ws.Range("M84").Value = -18857.667
$ws.Range("H113").Value = 898.95
$ws.Range("I113").Value = 798.9231
$ws.Range("J113").Value = 1084.7142
$ws.Range("K113").Value = 2396.7693
$ws.Range("L113").Value = 3254.1426
$ws.Range("M113").Value = -226.7692999999999
$ws.Range("N113").Value = -7594.142599999999
$ws.Range("H126").Value = 2148.4
$ws.Range("I126").Value = 2278.318
$ws.Range("K126").Value = 6834.954000000001
$ws.Range("M126").Value = -4364.954000000001

